$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.601.93"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.035.27"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.599"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.372"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0751"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "2.339.22"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").Value = "2.039.04"
$ws.Range("E18").Value = "  +3.23%  "
$ws.Range("D19").Value = "36.783.51"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  +18.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "0.0₃0796"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.126"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.30%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  +6.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +36.15%  "
$ws.Range("D43").Value = "1.477.23"
$ws.Range("E43").Value = "  +4.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0945"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.60%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.51%  "
